$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set the Notes (column F) values in the order that causes the shared-string
# table to be populated in the expected sequence:
#   44 Stringy
#   45 Stringy, but cleans up fine
#   46 Stringy, but cleans up pretty well
#   47 Not great quality
#   48 Quite blobby, cleans up well
$ws.Range("F15").Value = "Stringy"
$ws.Range("F7").Value = "Stringy, but cleans up fine"
$ws.Range("F25").Value = "Stringy, but cleans up pretty well"
$ws.Range("F20").Value = "Not great quality"
$ws.Range("F16").Value = "Quite blobby, cleans up well"
$ws.Range("F24").Value = "Not great quality"

# Update the selected cell to match the saved view state.
$ws.Range("F17").Select()
